# Applies: "adding averages and more checks"
# - Title & header-row fonts become bold + white (shared single font, losing the
#   old distinct "title sz14" font) on both sheets
# - Training Dashboard: PERIOD TO EXPIRE (H) and LAST UPDATE (I) refreshed to a
#   new reference date (16-Sep-2025), shifting PERIOD TO EXPIRE down by 8 days
# - Exam Dashboard: COMMENTS (E) column narrowed and messages updated to
#   "date is valid"; column E width reduced

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Training Dashboard")
$ws2 = $wb.Worksheets.Item("Exam Dashboard")

# ---------------------------------------------------------------------------
# 1. Styling: title (A1) and header row (row 2) fonts -> bold + white, same font
# ---------------------------------------------------------------------------
foreach ($ws in @($ws1, $ws2)) {
    $ws.Range("A1").Font.Size = 11
    $ws.Range("A1").Font.Color = 16777215
}

# Training Dashboard header spans A:K, Exam Dashboard header spans A:G
$ws1.Range("A2:K2").Font.Color = 16777215
$ws2.Range("A2:G2").Font.Color = 16777215

# ---------------------------------------------------------------------------
# 2. Training Dashboard: update PERIOD TO EXPIRE (H) and LAST UPDATE (I)
# ---------------------------------------------------------------------------
$newPeriod = @{
    3=399; 4=520; 5=495; 6=521; 7=519; 8=522; 9=520; 10=495; 11=520; 12=521;
    13=521; 14=524; 15=489; 16=523; 17=524; 18=522; 19=595; 20=584; 21=595;
    22=602; 23=678; 24=632; 25=631; 26=583; 27=665; 28=639; 29=636; 30=646;
    31=664; 32=637; 33=666; 34=594; 35=567; 36=679; 37=82; 38=-102; 39=-362;
    40=-49; 41=175; 42=175; 43=235; 44=232
}

for ($row = 3; $row -le 44; $row++) {
    $ws1.Cells.Item($row, 8).Value = $newPeriod[$row]
    # Leading apostrophe forces the date-looking text to stay plain text
    # (matching the workbook's existing inline-string date cells) instead of
    # being auto-converted into a date serial number.
    $ws1.Cells.Item($row, 9).Value = "'16-Sep-2025"
}

# ---------------------------------------------------------------------------
# 3. Exam Dashboard: narrow COMMENTS column and refresh comments
# ---------------------------------------------------------------------------
$ws2.Columns.Item(5).ColumnWidth = 14.17

$ws2.Cells.Item(3, 5).Value = "date is valid"
$ws2.Cells.Item(4, 5).Value = "date is valid"
$ws2.Cells.Item(5, 5).Value = "date is valid"
$ws2.Cells.Item(6, 5).Value = "date is valid"
$ws2.Cells.Item(7, 5).Value = "date is valid"
